$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Title paragraph: merge the three runs "Blatt 5 " + "Aufgabe" + " 2:"
#    into a single run "Blatt 5 Aufgabe 2:" (also drops the spell-check
#    proofErr markers automatically, since Find/Replace re-serialises the
#    matched range as one run).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Blatt 5 Aufgabe 2:", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Blatt 5 Aufgabe 2:", 2) | Out-Null

# Remove the _GoBack bookmark that used to sit at the end of the title
# paragraph - it gets re-added further down, near the end of the document.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2) Collapse the "Chrome developer tools" run sequence (and its proofErr
#    spell-check wrappers) into a single run.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute( `
    "Um die Lade- und Darstellungszeiten zu messen bieten sich insbesondere die Chrome developer tools an. ", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Um die Lade- und Darstellungszeiten zu messen bieten sich insbesondere die Chrome developer tools an. ", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Collapse "Des weiteren" run sequence into a single run.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute( `
    "Dabei wird genau aufgezeigt welche Datei wie lange benötigt und wie groß sie sind. Des weiteren ist gut erkennbar in welcher Reihenfolge die Dateien geladen werden.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Dabei wird genau aufgezeigt welche Datei wie lange benötigt und wie groß sie sind. Des weiteren ist gut erkennbar in welcher Reihenfolge die Dateien geladen werden.", 2) | Out-Null

# ---------------------------------------------------------------------------
# 4) Collapse "html code mit hervorhebung" run sequence into a single run.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute( `
    "Eine Einsicht in den html code mit hervorhebung der jeweiligen Elemente bietet der DOM Explorer.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Eine Einsicht in den html code mit hervorhebung der jeweiligen Elemente bietet der DOM Explorer.", 2) | Out-Null

# ---------------------------------------------------------------------------
# 5) Collapse "javascript" run sequence into a single run.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute( `
    "Weitere unter Umständen nützliche Tools existieren für die Echtzeitbearbeitung von javascript Dateien, die Emulation von anderen Browsern, Bildschirmen und Darstellungen und auch eine Konsolenausgabe.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Weitere unter Umständen nützliche Tools existieren für die Echtzeitbearbeitung von javascript Dateien, die Emulation von anderen Browsern, Bildschirmen und Darstellungen und auch eine Konsolenausgabe.", 2) | Out-Null

# ---------------------------------------------------------------------------
# 6) Append, after the last paragraph:
#      - two empty paragraphs
#      - a paragraph holding the (re-added) _GoBack bookmark
#      - a new bold paragraph about Apache JMeter
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$r = $lastPara.Range
$r.Collapse(0)
$r.InsertParagraphAfter()

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$r = $lastPara.Range
$r.Collapse(0)
$r.InsertParagraphAfter()

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$r = $lastPara.Range
$r.Collapse(0)
$r.InsertParagraphAfter()

$bookmarkPara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$d.Bookmarks.Add("_GoBack", $bookmarkPara.Range) | Out-Null

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newRange = $lastPara.Range
$newRange.Text = "Für Lasttests der Server  bietet sich Apache JMeter an, mit welchem man relativ einfach mehrere User die auf das System zugreifen simulieren kann. "
$newRange.Font.Bold = $true
